$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous contents (old range extended to AD19) before writing the new layout.
$ws.UsedRange.Clear()

# Target data for A1:T23 (row 1 = column index header, column A = row index, column B = row label).
# Each line is one worksheet row (1..23); cells are tab-separated.
# Token prefix: "E" = leave empty, "N:" = numeric value, "S:" = text (shared-string) value.
$table = @"
E	N:0	N:1	N:2	N:3	N:4	N:5	N:6	N:7	N:8	N:9	N:10	N:11	N:12	N:13	N:14	N:15	N:16	N:17	N:18
N:0	S:HKL	S:[2, 2, 2]	S:[3, 1, 0]	S:[1, 1, 0]	S:[3, 2, 1]	S:[4, 0, 0]	S:[2, 1, 1]	S:[2, 2, 0]	S:[2, 0, 0]	S:1Pair-A	S:1Pair-B	S:2Pairs-A	S:2Pairs-B	S:3Pairs-A	S:3Pairs-B	S:3Pairs-C	S:4Pairs	S:5A4F	S:MaxUnique
N:1	S:BT8Hex_2.5	N:0.999563562022222	N:1.000359118065113	N:0.9998377576235938	N:0.9998368180394401	N:1.000652386336827	N:0.9998364315635606	N:0.9998377576235938	N:1.000652386336827	N:0.9998377576235938	N:0.9998364315635606	N:1.000244408950194	N:1.000244408950194	N:1.000282645321833	N:1.000108858507994	N:1.000108858507994	N:1.000041083286894	N:1.000041083286894	N:1.000014345608459
N:2	S:BT8Hex_5	N:0.9991572253540326	N:1.000693266117362	N:0.9996870061908255	N:0.9996849935006176	N:1.001259290310473	N:0.9996841637442594	N:0.9996870061908255	N:1.001259290310473	N:0.9996870061908255	N:0.9996841637442594	N:1.000471727027366	N:1.000471727027366	N:1.000545573390698	N:1.000210153415186	N:1.000210153415186	N:1.000079366609096	N:1.000079366609096	N:1.000027657536261
N:3	S:BT8Hex_10	N:0.9983856444317494	N:1.001328584562606	N:0.9993995255154523	N:0.9993962805780844	N:1.002413678930295	N:0.9993949433916098	N:0.9993995255154523	N:1.002413678930295	N:0.9993995255154523	N:0.9993949433916098	N:1.000904311160952	N:1.000904311160952	N:1.00104573562817	N:1.000402715945786	N:1.000402715945786	N:1.000151918338202	N:1.000151918338202	N:1.000053109568299
N:4	S:BT8Hex_15	N:0.9976325414831287	N:1.001948618995555	N:0.9991190332001251	N:0.9991145246853895	N:1.003540256824746	N:0.9991126602320963	N:0.9991190332001251	N:1.003540256824746	N:0.9991190332001251	N:0.9991126602320963	N:1.001326458528421	N:1.001326458528421	N:1.001533845350799	N:1.000590650085656	N:1.000590650085656	N:1.000222745864273	N:1.000222745864273	N:1.00007793923684
N:5	S:Spiral2.5	N:0.9999749007722888	N:1.000020660747454	N:0.9999906583001333	N:0.9999906120094975	N:1.000037540387907	N:0.9999905934030382	N:0.9999906583001333	N:1.000037540387907	N:0.9999906583001333	N:0.9999905934030382	N:1.000014066895472	N:1.000014066895472	N:1.000016264846133	N:1.00000626403036	N:1.00000626403036	N:1.000002362597803	N:1.000002362597803	N:1.000000827603386
N:6	S:Spiral5	N:0.9999561914675327	N:1.000035789715006	N:0.9999841035433287	N:0.999983752197757	N:1.00006486637573	N:0.9999836080711139	N:0.9999841035433287	N:1.00006486637573	N:0.9999841035433287	N:0.9999836080711139	N:1.000024237223422	N:1.000024237223422	N:1.00002808805395	N:1.000010859330058	N:1.000010859330058	N:1.000004170383375	N:1.000004170383375	N:1.000001385228411
N:7	S:Spiral7.5	N:0.9999432552498964	N:1.000045841980631	N:0.9999801862454137	N:0.9999792176589232	N:1.000082776187046	N:0.9999788194933946	N:0.9999801862454137	N:1.000082776187046	N:0.9999801862454137	N:0.9999788194933946	N:1.00003079784022	N:1.00003079784022	N:1.000035812553691	N:1.000013927308618	N:1.000013927308618	N:1.000005492042817	N:1.000005492042817	N:1.000001682802551
N:8	S:Spiral10	N:0.999872131745729	N:1.000103970017084	N:0.9999543408891702	N:0.9999528253342344	N:1.000188140411769	N:0.9999522030446029	N:0.9999543408891702	N:1.000188140411769	N:0.9999543408891702	N:0.9999522030446029	N:1.000070171728186	N:1.000070171728186	N:1.000081437824485	N:1.000031561448514	N:1.000031561448514	N:1.000012256308678	N:1.000012256308678	N:1.000003935240432
N:9	S:Spiral15	N:0.9998102698548772	N:1.00015238508475	N:0.9999350865478505	N:0.9999309641166493	N:1.00027462143927	N:0.9999292682046509	N:0.9999350865478505	N:1.00027462143927	N:0.9999350865478505	N:0.9999292682046509	N:1.000101944821961	N:1.000101944821961	N:1.00011875824289	N:1.000046325397257	N:1.000046325397257	N:1.000018515684905	N:1.000018515684905	N:1.000005432541341
N:10	S:OffsetF45	N:0.9982055897701504	N:1.001426698866869	N:0.9994079411692615	N:0.9993544825523902	N:1.002562255068849	N:0.9993324848991358	N:0.9994079411692615	N:1.002562255068849	N:0.9994079411692615	N:0.9993324848991358	N:1.000947369983993	N:1.000947369983993	N:1.001107146278285	N:1.000434227045749	N:1.000434227045749	N:1.000177655576627	N:1.000177655576627	N:1.000048242054443
N:11	S:OffsetA45	N:0.9995672249145106	N:1.000379627028692	N:0.9998036936522696	N:0.9998261842488306	N:1.000703587174933	N:0.999835442414725	N:0.9998036936522696	N:1.000703587174933	N:0.9998036936522696	N:0.999835442414725	N:1.000269514794829	N:1.000269514794829	N:1.000306218872783	N:1.000114241080643	N:1.000114241080643	N:1.000036604223549	N:1.000036604223549	N:1.000019293238994
N:12	S:OffsetFTD	N:0.9896502829154836	N:1.008898876915657	N:0.9955761297980881	N:0.9959349935923743	N:1.016392937007566	N:0.9960826735632273	N:0.9955761297980881	N:1.016392937007566	N:0.9955761297980881	N:0.9960826735632273	N:1.006237805285397	N:1.006237805285397	N:1.00712482916215	N:1.002683913456294	N:1.002683913456294	N:1.000906967541742	N:1.000906967541742	N:1.000422648965399
N:13	S:OffsetATD	N:0.9968784428149133	N:1.002427378172071	N:0.9990521781334182	N:0.9989048862401907	N:1.00432592511089	N:0.9988442654986501	N:0.9990521781334182	N:1.00432592511089	N:0.9990521781334182	N:0.9988442654986501	N:1.00158509530477	N:1.00158509530477	N:1.001865856260537	N:1.000740789580986	N:1.000740789580986	N:1.000318636719094	N:1.000318636719094	N:1.000072179328356
N:14	S:Holden2.5	N:0.9914949215434024	N:1.007001490445431	N:0.9968334892827655	N:0.9968183557812798	N:1.01272099135825	N:0.9968121218380639	N:0.9968334892827655	N:1.01272099135825	N:0.9968334892827655	N:0.9968121218380639	N:1.004766556598157	N:1.004766556598157	N:1.005511534547248	N:1.00212220082636	N:1.00212220082636	N:1.000800022940461	N:1.000800022940461	N:1.000280228374865
N:15	S:Holden5	N:0.993058030690656	N:1.005719080207106	N:0.9974088785223879	N:0.9974008693571582	N:1.010393570635947	N:0.9973975694744341	N:0.9974088785223879	N:1.010393570635947	N:0.9974088785223879	N:0.9973975694744341	N:1.003895570055191	N:1.003895570055191	N:1.004503406772496	N:1.001733339544256	N:1.001733339544256	N:1.000652224288789	N:1.000652224288789	N:1.000229666481282
N:16	S:Holden10	N:0.9962239132922237	N:1.003120895629267	N:0.9985754982560598	N:0.9985810992918772	N:1.005677685355994	N:0.9985834024100173	N:0.9985754982560598	N:1.005677685355994	N:0.9985754982560598	N:0.9985834024100173	N:1.002130543883006	N:1.002130543883006	N:1.002460661131759	N:1.000945528674024	N:1.000945528674024	N:1.000353021069533	N:1.000353021069533	N:1.000127082372573
N:17	S:Holden15	N:0.995973951997324	N:1.003332894085547	N:0.9984730555203349	N:0.9984844118287821	N:1.006066559994675	N:0.9984890864013287	N:0.9984730555203349	N:1.006066559994675	N:0.9984730555203349	N:0.9984890864013287	N:1.002277823198002	N:1.002277823198002	N:1.00262951349385	N:1.001009567305446	N:1.001009567305446	N:1.000375439359168	N:1.000375439359168	N:1.000136659971332
N:18	S:HexGrid-90degTilt2.5degRes	N:0.999999642087146	N:1.000001526876128	N:0.9999980120544331	N:0.9999992381471445	N:1.000003507066983	N:0.999999743083752	N:0.9999980120544331	N:1.000003507066983	N:0.9999980120544331	N:0.999999743083752	N:1.000001625075368	N:1.000001625075368	N:1.000001592342288	N:1.000000420735056	N:1.000000420735056	N:0.9999998185649003	N:0.9999998185649003	N:1.000000278219265
N:19	S:HexGrid-90degTilt5degRes	N:0.9999599147259965	N:1.000031345648932	N:0.9999875669404142	N:0.9999858489540983	N:1.000055974196975	N:0.9999851420782735	N:0.9999875669404142	N:1.000055974196975	N:0.9999875669404142	N:0.9999851420782735	N:1.000020558137624	N:1.000020558137624	N:1.000024153974727	N:1.000009561071888	N:1.000009561071888	N:1.000004062539019	N:1.000004062539019	N:1.000000965424115
N:20	S:HexGrid-90degTilt10degRes	N:0.9998697462216242	N:1.000110696069929	N:0.9999462815488652	N:0.9999495054421271	N:1.000203183820438	N:0.999950828013853	N:0.9999462815488652	N:1.000203183820438	N:0.9999462815488652	N:0.999950828013853	N:1.000077005917146	N:1.000077005917146	N:1.000088235968074	N:1.000033431127719	N:1.000033431127719	N:1.000011643733006	N:1.000011643733006	N:1.00000504018614
N:21	S:HexGrid-90degTilt15degRes	N:0.9997012196973705	N:1.000259525235631	N:0.9998683310941364	N:0.999881308274888	N:1.000479579573885	N:0.9998866459810432	N:0.9998683310941364	N:1.000479579573885	N:0.9998683310941364	N:0.9998866459810432	N:1.000183112777464	N:1.000183112777464	N:1.000208583596853	N:1.000078185549688	N:1.000078185549689	N:1.0000257219358	N:1.0000257219358	N:1.000012768309492
"@

$lines = $table -split "`n"

for ($r = 0; $r -lt $lines.Length; $r++) {
    $line = $lines[$r].TrimEnd("`r")
    $cells = $line -split "`t"
    for ($c = 0; $c -lt $cells.Length; $c++) {
        $tok = $cells[$c]
        if ($tok -eq "E") {
            continue
        } elseif ($tok.StartsWith("N:")) {
            $val = [double]($tok.Substring(2))
            $ws.Cells.Item($r + 1, $c + 1).Value = $val
        } elseif ($tok.StartsWith("S:")) {
            $val = $tok.Substring(2)
            $ws.Cells.Item($r + 1, $c + 1).Value = $val
        }
    }
}

# Re-apply the header style (bold font, thin border, centered/top-aligned) that Excel
# used for the index row (row 1, columns B:T) and the row-index column (A2:A23).
$headerRow = $ws.Range("B1:T1")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160
$headerRow.Borders.LineStyle = 1
$headerRow.Borders.Weight = 2

$indexCol = $ws.Range("A2:A23")
$indexCol.Font.Bold = $true
$indexCol.HorizontalAlignment = -4108
$indexCol.VerticalAlignment = -4160
$indexCol.Borders.LineStyle = 1
$indexCol.Borders.Weight = 2
